$d = $word.ActiveDocument
Write-Host "Paragraph count: $($d.Paragraphs.Count)"
